# Adds two new weekly price records ("Ají" / Vega Monumental Concepción) to the
# data table, inserted at row 136 (pushing the existing records, previously at
# rows 136-178, down to rows 138-180).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 136, shifting existing rows 136-178 down to 138-180
$ws.Rows("136:137").Insert()

# Populate the newly inserted row 136 with new data
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 45009
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100112021
$ws.Range("G136").Value = "Ají"
$ws.Range("H136").Value = "Americana (o)"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 100
$ws.Range("K136").Value = 21000
$ws.Range("L136").Value = 22000
$ws.Range("M136").Value = 21500
$ws.Range("N136").Value = "`$/saco 25 kilos"
$ws.Range("O136").Value = "Región Metropolitana"
$ws.Range("P136").Value = 860
$ws.Range("Q136").Value = 25
$ws.Range("R136").Value = "Hortaliza"
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 45009
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112021
$ws.Range("G137").Value = "Ají"
$ws.Range("H137").Value = "Americana (o)"
$ws.Range("I137").Value = "Segunda"
$ws.Range("J137").Value = 50
$ws.Range("K137").Value = 17000
$ws.Range("L137").Value = 17000
$ws.Range("M137").Value = 17000
$ws.Range("N137").Value = "`$/saco 25 kilos"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 680
$ws.Range("Q137").Value = 25
$ws.Range("R137").Value = "Hortaliza"
